# QR-13AF_data.xlsx - "New data and plots added"
#
# The dataset's last two columns (D = "off", E = "on"-type flag) are
# flipped for every data row: column D goes from 0 -> 1 and column E
# goes from 1 -> 0, for rows 1 through 68. The rest of the data
# (columns A-C) is untouched.
#
# The author also scrolled/selected a different cell (F55) before
# saving, which is reflected in the worksheet's selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flip the D/E indicator columns for all 68 data rows in one shot.
$ws.Range("D1:D68").Value = 1
$ws.Range("E1:E68").Value = 0

# Update the active selection to match where the author left off.
$ws.Range("F55").Select()
